$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at J:K (old J..N shift right to L..P)
$ws.Range("J1:K1").EntireColumn.Insert()

# New column headers / values (J, K are brand new columns)
$ws.Range("J1").Value = "Avai. C fluxes from"
$ws.Range("K1").Value = "PhenoCam. measure."
$ws.Range("J2").Value = "1996-"
$ws.Range("K2").Value = "yes"
$ws.Range("J3").Value = "1996-"
$ws.Range("J4").Value = "1996-"
$ws.Range("K4").Value = "yes"
$ws.Range("J5").Value = "1997-"
$ws.Range("K5").Value = "yes"

# Updated header text for the (shifted) "other available measurements" column
$ws.Range("N1").Value = "Other avai. Measure."

# Updated data values
$ws.Range("N2").Value = "sapflow, dendrometer "
$ws.Range("O2").Value = "Timo Vesala, Albert Porcar-Castell, Anna Lintunen"
$ws.Range("N4").Value = "dendrometer, litterfall"

# Hyperlinks: the engine keeps hyperlink *ranges* pinned to their original
# cells rather than shifting them with EntireColumn.Insert, so re-create them
# pointing at the new location (P2:P4).
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("P2"), "mailto:joan.porcar@helsinki.fi;timo.vesala@helsinki.fi", [System.Type]::Missing, [System.Type]::Missing, "joan.porcar@helsinki.fi;timo.vesala@helsinki.fi     ")
$ws.Hyperlinks.Add($ws.Range("P3"), "mailto:michiel.vandermolen@wur.nl")
$ws.Hyperlinks.Add($ws.Range("P4"), "mailto:thomas.gruenwald@tu-dresden.de")

$ws.Columns.AutoFit()

$ws.Range("N7").Select()
